# Apply updated dSF (column F) values for the specified rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 0
    7  = -8
    8  = -2
    18 = -4
    21 = -3
    22 = -2
    27 = 6
    28 = -2
    34 = 8
    40 = -2
    43 = 6
    44 = -4
    45 = -2
    47 = 3
    49 = 8
    57 = 5
    63 = -1
    64 = -10
    65 = -1
    66 = -3
    69 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
